# Apply the data edit described by the diff: recode job-title columns
# (C, D, E), the numeric column F, and the ordering/ranking column G
# for rows 1-20 of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Cashier"
$ws.Range("D1").Value = "Stock Clerk"
$ws.Range("E1").Value = "Customer Service"
$ws.Range("F1").Value = 9
$ws.Range("G1").Value = "2,4,6,1,5,3"

$ws.Range("C2").Value = "Customer Service"
$ws.Range("D2").Value = "Stock Clerk"
$ws.Range("E2").Value = "Cashier"
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = "3,1,5,2"

$ws.Range("C3").Value = "Stock Clerk"
$ws.Range("D3").Value = "Customer Service"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 18
$ws.Range("G3").Value = "6,3,1,2,5"

$ws.Range("C4").Value = "Cashier"
$ws.Range("D4").Value = "Customer Service"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 17
$ws.Range("G4").Value = "6,5,4,2,3,1"

$ws.Range("C5").Value = "Cashier"
$ws.Range("D5").Value = "Stock Clerk"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 17
$ws.Range("G5").Value = "3,2,4,5,1,6"

$ws.Range("C6").Value = "Cashier"
$ws.Range("D6").Value = "Stock Clerk"
$ws.Range("E6").Value = "Customer Service"
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = "1,4,2,3"

$ws.Range("C7").Value = "Cashier"
$ws.Range("D7").Value = "Stock Clerk"
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = 20
$ws.Range("G7").Value = "5,6,4,3,2,1"

$ws.Range("C8").Value = "Cashier"
$ws.Range("D8").Value = "Customer Service"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = 18
$ws.Range("G8").Value = "4,6,2,3,1,5"

$ws.Range("C9").Value = "Stock Clerk"
$ws.Range("D9").Value = "Cashier"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = 11
$ws.Range("G9").Value = "3,2,6,4"

$ws.Range("C10").Value = "Cashier"
$ws.Range("D10").Value = "Stock Clerk"
$ws.Range("E10").Value = "Customer Service"
$ws.Range("F10").Value = 12
$ws.Range("G10").Value = "4,5,3,2"

$ws.Range("C11").Value = "Cashier"
$ws.Range("D11").Value = "Customer Service"
$ws.Range("E11").Value = "Stock Clerk"
$ws.Range("F11").Value = 19
$ws.Range("G11").Value = "2,6,5,1,4,3"

$ws.Range("C12").Value = "Cashier"
$ws.Range("D12").Value = "Customer Service"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 11
$ws.Range("G12").Value = "4,1,3,2"

$ws.Range("C13").Value = "Stock Clerk"
$ws.Range("D13").Value = "Customer Service"
$ws.Range("E13").Value = "Cashier"
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = "6,2,1,5,4,3"

$ws.Range("C14").Value = "Cashier"
$ws.Range("D14").Value = "Stock Clerk"
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = "3,2,6,1"

$ws.Range("C15").Value = "Stock Clerk"
$ws.Range("D15").Value = "Cashier"
$ws.Range("E15").Value = "Customer Service"
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = "3,1,6,5,4,2"

$ws.Range("C16").Value = "Customer Service"
$ws.Range("D16").Value = "Stock Clerk"
$ws.Range("E16").Value = "Cashier"
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = "5,1,4,2,3"

$ws.Range("C17").Value = "Customer Service"
$ws.Range("D17").Value = "Cashier"
$ws.Range("E17").Value = "Stock Clerk"
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = "1,6,3,2,4,5"

$ws.Range("C18").Value = "Cashier"
$ws.Range("D18").Value = "Stock Clerk"
$ws.Range("E18").Value = "Customer Service"
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = "3,2,6,1,5,4"

$ws.Range("C19").Value = "Customer Service"
$ws.Range("D19").Value = "Cashier"
$ws.Range("E19").Value = "Stock Clerk"
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = "2,6,1,3,4"

$ws.Range("C20").Value = "Stock Clerk"
$ws.Range("D20").Value = "Customer Service"
$ws.Range("E20").Value = "Cashier"
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = "5,4,1,2"

